# Rebuild the invoices header row with the new "basic module" schema and
# drop the two sample data rows, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the existing bold/bordered header style (currently used by
#        A1:H1, style index s="1") onto the columns that will be newly
#        added (I1:S1) before we touch their values, so every header cell
#        in A1:S1 ends up sharing the same cell style. ---
$ws.Range("A1").Copy()
$ws.Range("I1:S1").PasteSpecial(-4122)   # xlPasteFormats

# --- 2. Write the full new header row, A1:S1. ---
$ws.Range("A1").Value  = "id"
$ws.Range("B1").Value  = "invoice_number"
$ws.Range("C1").Value  = "reference_type"
$ws.Range("D1").Value  = "reference_id"
$ws.Range("E1").Value  = "date"
$ws.Range("F1").Value  = "due_date"
$ws.Range("G1").Value  = "customer"
$ws.Range("H1").Value  = "billing_address"
$ws.Range("I1").Value  = "shipping_address"
$ws.Range("J1").Value  = "items"
$ws.Range("K1").Value  = "subtotal"
$ws.Range("L1").Value  = "discount"
$ws.Range("M1").Value  = "gst_breakdown"
$ws.Range("N1").Value  = "total_gst"
$ws.Range("O1").Value  = "total_amount"
$ws.Range("P1").Value  = "payment_status"
$ws.Range("Q1").Value  = "notes"
$ws.Range("R1").Value  = "created_at"
$ws.Range("S1").Value  = "updated_at"

# --- 3. Remove the two sample data rows (2 and 3) entirely. ---
$ws.Range("A2:H3").ClearContents()
